$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.4092363652084714,
    -0.00999999930138884,
    -0.087810384582763135,
    -0.011999999815550666,
    -0.0059999992885888531,
    -0.0059999992683970049,
    -0.019999999147930936,
    -0.019999999144624248,
    -0.005999999261522504,
    -0.0059999992603323449,
    -0.0044999992731433736,
    -0.0059999992605739294,
    -0.0059999992624657494,
    0.0072298850662484426,
    0.025969758352838213,
    -0.0059999992620407561,
    -0.005999999259337585,
    -0.0089999992333220646,
    -0.0089999993094300734,
    -0.0089999993032963133,
    -0.008999999302379269,
    -0.0089999993018241575,
    -0.0089999992751019775,
    -0.041999998983234477,
    -0.041999998977789943,
    -0.0059999992666064372,
    -0.0059999992638486432,
    -0.005999999249814536,
    -0.011999999189470145,
    -0.01944568505171107,
    -0.014999999152339072,
    -0.020999999100417277,
    -0.0059999992277575132
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
